# The author's edit reordered two pairs of shared-string table entries:
#   - "Northwest reels from deadly back-to-back storms" and
#     "Washington and Oregon Residents Fight to Save Homes From Muddy Waters" swapped places
#   - the CNN and Fox News article URLs swapped places
# Because the worksheet cells keep referencing the same shared-string slots,
# the net, observable effect is that the text shown in A4/A5 (titles) and
# E4/E5 (uris) is swapped, while the hyperlink targets (row-bound relationship
# ids) are left exactly as they were. Reproduce that here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titleA4 = $ws.Range("A4").Value2
$titleA5 = $ws.Range("A5").Value2
$uriE4 = $ws.Range("E4").Value2
$uriE5 = $ws.Range("E5").Value2

$ws.Range("A4").Value = $titleA5
$ws.Range("A5").Value = $titleA4

$ws.Range("E4").Value = $uriE5
$ws.Range("E5").Value = $uriE4
